$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds the date as literal text (e.g. "2025/10/06"), not a real
# Excel date serial. Format the cell as Text before typing the value so the
# engine doesn't auto-convert the "yyyy/mm/dd"-looking string into a date,
# then clear the formatting back to the sheet's default (no explicit style)
# so the new row matches the plain, unstyled data rows above it.
$ws.Range("A66").NumberFormat = "@"
$ws.Range("A66").Value = "2025/10/06"
$ws.Range("A66").ClearFormats()

$ws.Range("B66").Value = "月"
$ws.Range("C66").Value = 1
$ws.Range("D66").Value = 65
